$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.623.64"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "3.380.33"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.09"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.99"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "3.376.35"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.632"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.42"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -3.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.20"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "3.915.45"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.396.06"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.119"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.18"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "65.519.69"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.83"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.996"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.67"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.88"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -4.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.75"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.29"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +6.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.09"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.91"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.59"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.70"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.12"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.56"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.41"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "576.28"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.93"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.108"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.141"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.78"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.375"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "0.0₃0737"
$ws.Range("E41").Value = "  -3.22%  "
$ws.Range("D42").Value = "3.107.69"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.81"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0416"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.134"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.16"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.43"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -2.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.45"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.47"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +1.27%  "
